# Appends the new 2026/01/01 - 2026/01/05 "sei2" rows (545-570) to Sheet1,
# extending the sheet's data range from A1:D544 to A1:D570.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstNewRow = 545

$newRows = @(
    @("2026/01/01", "木", 2,  201),
    @("2026/01/01", "木", 5,  201),
    @("2026/01/01", "木", 13, 201),
    @("2026/01/01", "木", 16, 201),
    @("2026/01/01", "木", 19, 201),
    @("2026/01/02", "金", 1,  201),
    @("2026/01/02", "金", 5,  201),
    @("2026/01/02", "金", 8,  201),
    @("2026/01/02", "金", 13, 201),
    @("2026/01/02", "金", 16, 201),
    @("2026/01/02", "金", 19, 201),
    @("2026/01/02", "金", 22, 201),
    @("2026/01/03", "土", 1,  201),
    @("2026/01/03", "土", 4,  201),
    @("2026/01/03", "土", 7,  201),
    @("2026/01/03", "土", 13, 201),
    @("2026/01/03", "土", 16, 201),
    @("2026/01/03", "土", 19, 201),
    @("2026/01/03", "土", 22, 201),
    @("2026/01/04", "日", 2,  201),
    @("2026/01/04", "日", 4,  201),
    @("2026/01/04", "日", 7,  201),
    @("2026/01/04", "日", 13, 201),
    @("2026/01/04", "日", 22, 201),
    @("2026/01/05", "月", 2,  201),
    @("2026/01/05", "月", 7,  201)
)

$lastNewRow = $firstNewRow + $newRows.Count - 1

# Column A holds date-like text (e.g. "2026/01/01"). A plain .Value assignment
# of a date-shaped string gets auto-coerced by Excel into a date serial
# number, but the source data stores these as literal text. Forcing the
# "Text" number format before assignment keeps the literal string; the style
# is reset back to "Normal" afterwards so the cells end up with no explicit
# style index, same as every other data row already on the sheet.
$dateRange = $ws.Range("A$firstNewRow" + ":A$lastNewRow")
$dateRange.NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $firstNewRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$dateRange.Style = "Normal"

Write-Host "UsedRange: $($ws.UsedRange.Address())"
